# Regen save_data: update column G ("K") values for rows 2-27 with the
# newly-calculated strikeout (K) counts, replacing the old Strike# values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new K value (column G)
$kValues = @{
    2  = 1
    3  = 2
    4  = 0
    5  = 0
    6  = 1
    7  = 0
    8  = 1
    9  = 0
    10 = 1
    11 = 0
    12 = 1
    13 = 1
    14 = 0
    15 = 1
    16 = 2
    17 = 0
    18 = 0
    19 = 0
    20 = 1
    21 = 2
    22 = 1
    23 = 1
    24 = 0
    25 = 1
    26 = 3
    27 = 1
}

foreach ($row in $kValues.Keys | Sort-Object) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
